# Ajout de l'activité du jour : nouvelle ligne 60 dans le journal de travail.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reprend la mise en forme de la dernière ligne existante (59) pour la
# nouvelle ligne 60 (style date colonne A, style texte avec retour à la
# ligne colonne B).
$ws.Range("A59:C59").Copy()
$ws.Range("A60:C60").PasteSpecial(-4122)

# Date de l'activité (29/03/2018)
$ws.Range("A60").Value = 43188

# Description de l'activité du jour
$ws.Range("B60").Value = "Ajout de différentes informations sur la page d'acceuil du site, j'ai ajouté les trois catégories que je propose sur le site comme ça les gens n'ont plus qu'à cliquer sur la catégorie qui les intéresse et seront dirgié vers celle-ci. Ensuite en bas de page j'ai affiché une dizaine d'article que je récupère aléatoirement dans la base de données. Les gens clique sur le bouton et sont dirigés vers l'article pour le commander."

# Durée de l'activité (même texte que les autres lignes "4 périodes")
$ws.Range("C60").Value = "4 périodes"

# Hauteur de la ligne, comme dans le classeur final
$ws.Rows.Item(60).RowHeight = 60

# Zoom de la feuille
$excel.ActiveWindow.Zoom = 143

# La sélection se retrouve sur la cellule suivante, prête pour la prochaine entrée
$ws.Range("C61").Select()

Write-Host "Ligne 60 ajoutée"
